$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting so values
# like "8.50" or "0.0508" are not coerced into numbers, losing trailing
# zeros / dot-grouping (e.g. "59.394.79").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.394.79"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "3.316.05"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "408.70"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "111.41"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "0.585"
$ws.Range("E7").Value = "  +4.62%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "3.838.86"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "8.50"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "19.22"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "3.318.53"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "59.178.30"
$ws.Range("E18").Value = "  +3.66%  "
$ws.Range("D19").Value = "10.70"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("D20").Value = "3.33"
$ws.Range("E20").Value = "  +0.59%  "
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").Value = "13.05"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "301.84"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "75.29"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "28.43"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").Value = "0.180"
$ws.Range("E28").Value = "  +5.76%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.72"
$ws.Range("E29").Value = "  +26.82%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "11.56"
$ws.Range("E34").Value = "  +4.88%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "39.67"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.0508"
$ws.Range("E36").Value = "  +4.96%  "
$ws.Range("D37").Value = "51.87"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("D41").Value = "138.53"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.283"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.91"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "16.76"
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  +8.20%  "
$ws.Range("D48").Value = "22.22"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "2.202.02"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  -0.69%  "
